$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtered save games) - update B,C,D,E,G columns for rows 2-5
$ws.Range("B2").Value = 0.6753301551942219
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 2.997429241610044

$ws.Range("B3").Value = 3.230985683306322
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 0.1575252929769615
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 13.71653804550039

$ws.Range("B4").Value = 0.6753301551942219
$ws.Range("C4").Value = 1.667794583268128
$ws.Range("D4").Value = 26.21740644021617
$ws.Range("E4").Value = 0.496779210170732
$ws.Range("G4").Value = 29.05731038884925

$ws.Range("B5").Value = 1.459612070389937
$ws.Range("C5").Value = 1.667794583268128
$ws.Range("D5").Value = 3.900430680208489
$ws.Range("E5").Value = 0.496779210170732
$ws.Range("G5").Value = 7.524616544037286
